# Clear the two paragraphs that contained the gradient-descent error values,
# leaving them as empty paragraphs (<w:p/>), as described in the commit
# "Eliminado error de los apartados de gradiente".

$d = $word.ActiveDocument

function Clear-ParagraphByText($marker) {
    $found = $false
    foreach ($p in $d.Paragraphs) {
        $text = $p.Range.Text
        if ($text -like "*$marker*") {
            $r = $p.Range
            # Trim off the trailing paragraph mark so we only delete the
            # paragraph's own content, leaving the paragraph itself (empty).
            $r.MoveEnd(1, -1) | Out-Null
            $r.Text = ""
            $found = $true
            break
        }
    }
    return $found
}

Clear-ParagraphByText("5.7833") | Out-Null
Clear-ParagraphByText("8.5838") | Out-Null
